# Revert "Lentekhi Municipality area" sheet from the 3-year (1989/2002/2014)
# layout back to the single-year (2014) layout, matching the target diff:
#  - drop the "(according to the population census data)" caption row text
#  - drop the 1989/2002 columns (C, D), keep only a 2014 column (B)
#  - tighten row heights to 20.1pt and extend sheetData down to row 9
#  - fix column widths back to A=15.71 / B=8.71

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Wipe all existing content/formatting and extra columns so we start
#    from a clean A:B sheet (old sheet used A:D).
# ---------------------------------------------------------------------
$ws.Cells.Clear()
$ws.Columns.Item(4).Delete()
$ws.Columns.Item(3).Delete()

$ws.Columns.Item(1).ColumnWidth = 15.7109375
$ws.Columns.Item(2).ColumnWidth = 8.7109375

# ---------------------------------------------------------------------
# 2. Values
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "Area of the municipality of Lentekhi"
$ws.Range("A3").Value = "(sq. km)"
$ws.Range("B4").Value = 2014
$ws.Range("A5").Value = "Area"
$ws.Range("B5").Value = 1344

# ---------------------------------------------------------------------
# 3. Row heights (rows 1-9, all 20.1pt custom height)
# ---------------------------------------------------------------------
for ($r = 1; $r -le 9; $r++) {
    $ws.Rows.Item($r).RowHeight = 20.1
}

# ---------------------------------------------------------------------
# 4. Formatting per-cell
# ---------------------------------------------------------------------

# --- A1: bold 11pt Arial title, left/vcenter, white fill ---
$a1 = $ws.Range("A1")
$a1.Interior.Pattern = 1
$a1.Interior.Color = 16777215
$a1.Font.Name = "Arial"
$a1.Font.Size = 11
$a1.Font.Bold = $true
$a1.HorizontalAlignment = -4131
$a1.VerticalAlignment = -4108
$a1.NumberFormat = "@"

# --- A2/B1/B2 (blank formatted row 1-2 cells use plain 10pt Arial, white fill) ---
$row1fmt = $ws.Range("B1")
$row1fmt.Interior.Pattern = 1
$row1fmt.Interior.Color = 16777215
$row1fmt.Font.Name = "Arial"
$row1fmt.Font.Size = 10
$row1fmt.Font.Bold = $false

$row2fmt = $ws.Range("A2:B2")
$row2fmt.Interior.Pattern = 1
$row2fmt.Interior.Color = 16777215
$row2fmt.Font.Name = "Arial"
$row2fmt.Font.Size = 10
$row2fmt.Font.Bold = $false

# --- A3: "(sq. km)" caption, left aligned, white fill, thin 10pt Arial ---
$a3 = $ws.Range("A3")
$a3.Interior.Pattern = 1
$a3.Interior.Color = 16777215
$a3.Font.Name = "Arial"
$a3.Font.Size = 10
$a3.Font.Bold = $false
$a3.HorizontalAlignment = -4131

# --- A4: empty cell with medium box border (left/right/top), white fill ---
$a4 = $ws.Range("A4")
$a4.Interior.Pattern = 1
$a4.Interior.Color = 16777215
$a4.Font.Name = "Arial"
$a4.Font.Size = 10
$a4.Font.Bold = $false
$a4.VerticalAlignment = -4108
$a4.Borders.Item(7).LineStyle = 1
$a4.Borders.Item(7).Weight = -4138
$a4.Borders.Item(10).LineStyle = 1
$a4.Borders.Item(10).Weight = -4138
$a4.Borders.Item(8).LineStyle = 1
$a4.Borders.Item(8).Weight = -4138
$a4.Borders.Item(9).LineStyle = -4142

# --- B4: "2014" header, centered, medium border box, white fill ---
$b4 = $ws.Range("B4")
$b4.Interior.Pattern = 1
$b4.Interior.Color = 16777215
$b4.Font.Name = "Arial"
$b4.Font.Size = 10
$b4.Font.Bold = $false
$b4.HorizontalAlignment = -4108
$b4.VerticalAlignment = -4108
$b4.Borders.Item(7).LineStyle = 1
$b4.Borders.Item(7).Weight = 2
$b4.Borders.Item(10).LineStyle = 1
$b4.Borders.Item(10).Weight = -4138
$b4.Borders.Item(8).LineStyle = 1
$b4.Borders.Item(8).Weight = -4138
$b4.Borders.Item(9).LineStyle = -4142

# --- A5: "Area" label, bold, vcenter, medium border (left/top/bottom), white fill ---
$a5 = $ws.Range("A5")
$a5.Interior.Pattern = 1
$a5.Interior.Color = 16777215
$a5.Font.Name = "Arial"
$a5.Font.Size = 10
$a5.Font.Bold = $true
$a5.VerticalAlignment = -4108
$a5.Borders.Item(7).LineStyle = 1
$a5.Borders.Item(7).Weight = -4138
$a5.Borders.Item(8).LineStyle = 1
$a5.Borders.Item(8).Weight = -4138
$a5.Borders.Item(9).LineStyle = 1
$a5.Borders.Item(9).Weight = -4138
$a5.Borders.Item(10).LineStyle = -4142

# --- B5: 1344 value, bold, "0.0" number format, right/vcenter, medium border (right/top/bottom) ---
$b5 = $ws.Range("B5")
$b5.Interior.Pattern = 1
$b5.Interior.Color = 16777215
$b5.Font.Name = "Arial"
$b5.Font.Size = 10
$b5.Font.Bold = $true
$b5.NumberFormat = "0.0"
$b5.HorizontalAlignment = -4152
$b5.VerticalAlignment = -4108
$b5.Borders.Item(10).LineStyle = 1
$b5.Borders.Item(10).Weight = -4138
$b5.Borders.Item(8).LineStyle = 1
$b5.Borders.Item(8).Weight = -4138
$b5.Borders.Item(9).LineStyle = 1
$b5.Borders.Item(9).Weight = -4138
$b5.Borders.Item(7).LineStyle = -4142

# --- rows 6-9: plain white-filled blank rows (default formatting, just height) ---

$ws.Range("A1").Select()
